$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $ws.Range("D6").Copy() | Out-Null
    $r.PasteSpecial(-4122) | Out-Null
}

$ws.Range("D2").Value = "26.901.14"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.546.74"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  +0.31%  "
Set-TextValue "D5" "206.51"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  -0.06%  "
Set-TextValue "D9" "21.34"
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("E10").Value = "  -0.21%  "
Set-TextValue "D11" "0.0858"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "1.766.56"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").Value = "1.549.57"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "26.886.91"
$ws.Range("E16").Value = "  -0.01%  "
Set-TextValue "D17" "61.41"
$ws.Range("E17").Value = "  +0.23%  "
Set-TextValue "D18" "214.78"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  +0.47%  "
Set-TextValue "D20" "7.20"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  -2.60%  "
Set-TextValue "D23" "9.19"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  -2.99%  "
Set-TextValue "D25" "151.77"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D33").Value = "1.356.42"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("E36").Value = "  +4.64%  "
Set-TextValue "D38" "0.0165"
$ws.Range("E38").Value = "  -0.09%  "
Set-TextValue "D39" "0.521"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("E42").Value = "  +3.94%  "
Set-TextValue "D43" "0.990"
$ws.Range("E43").Value = "  -0.79%  "
Set-TextValue "D44" "2.22"
$ws.Range("E44").Value = "  +2.07%  "
Set-TextValue "D45" "63.42"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("D47").Value = "1.681.36"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("E48").Value = "  -5.43%  "
Set-TextValue "D49" "85.79"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").Value = "0.0₇0966"
$ws.Range("E51").Value = "  -1.48%  "

$excel.CutCopyMode = $false
